# Insert two new price-report rows (Naranja / Macroferia Regional de Talca)
# right after the current row 528, shifting the existing rows 529:578 down to
# 531:580. The two new rows land at 529 and 530.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("529:530").Insert()

# --- New row 529 ---------------------------------------------------------
$ws.Cells.Item(529,1).Value  = 5
$ws.Cells.Item(529,2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(529,3).Value  = 'Maule'
$ws.Cells.Item(529,4).Value  = 44769
$ws.Cells.Item(529,5).Value  = 7
$ws.Cells.Item(529,6).Value  = 'Fruta'
$ws.Cells.Item(529,7).Value  = 100102
$ws.Cells.Item(529,8).Value  = 'Cítricos'
$ws.Cells.Item(529,9).Value  = 100102005
$ws.Cells.Item(529,10).Value = 'Naranja'
$ws.Cells.Item(529,11).Value = 'Fukumoto'
$ws.Cells.Item(529,12).Value = 'Primera'
$ws.Cells.Item(529,13).Value = 260
$ws.Cells.Item(529,14).Value = 6000
$ws.Cells.Item(529,15).Value = 6000
$ws.Cells.Item(529,16).Value = 6000
$ws.Cells.Item(529,17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(529,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(529,19).Value = 400
$ws.Cells.Item(529,20).Value = 15

# --- New row 530 ---------------------------------------------------------
$ws.Cells.Item(530,1).Value  = 5
$ws.Cells.Item(530,2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(530,3).Value  = 'Maule'
$ws.Cells.Item(530,4).Value  = 44769
$ws.Cells.Item(530,5).Value  = 7
$ws.Cells.Item(530,6).Value  = 'Fruta'
$ws.Cells.Item(530,7).Value  = 100102
$ws.Cells.Item(530,8).Value  = 'Cítricos'
$ws.Cells.Item(530,9).Value  = 100102005
$ws.Cells.Item(530,10).Value = 'Naranja'
$ws.Cells.Item(530,11).Value = 'Lane Late'
$ws.Cells.Item(530,12).Value = 'Primera'
$ws.Cells.Item(530,13).Value = 320
$ws.Cells.Item(530,14).Value = 5000
$ws.Cells.Item(530,15).Value = 5000
$ws.Cells.Item(530,16).Value = 5000
$ws.Cells.Item(530,17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(530,18).Value = 'Región de O''Higgins'
$ws.Cells.Item(530,19).Value = 333
$ws.Cells.Item(530,20).Value = 15
